$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cd9"
$ws.Cells.Item(2,3).Value = "L1cam"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 40.314886
$ws.Cells.Item(2,8).Value = 120.944658
$ws.Cells.Item(2,9).Value = 0.1963336494301312
$ws.Cells.Item(2,10).Value = 0.1963336494301312
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 11.88712433333333
$ws.Cells.Item(2,14).Value = 35.661373
$ws.Cells.Item(2,15).Value = 0.5967229292030898
$ws.Cells.Item(2,16).Value = 0.5967229292030898
$ws.Cells.Item(2,17).Value = 479.2280623661593
$ws.Cells.Item(2,18).Value = 4313.052561295433
$ws.Cells.Item(2,19).Value = 0.1171567903890804
$ws.Cells.Item(2,20).Value = 0.1171567903890804
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cd9"
$ws.Cells.Item(3,3).Value = "L1cam"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 40.314886
$ws.Cells.Item(3,8).Value = 120.944658
$ws.Cells.Item(3,9).Value = 0.1963336494301312
$ws.Cells.Item(3,10).Value = 0.1963336494301312
$ws.Cells.Item(3,11).Value = 2
$ws.Cells.Item(3,12).Value = 0.6666666666666666
$ws.Cells.Item(3,13).Value = 0.3785896666666667
$ws.Cells.Item(3,14).Value = 1.135769
$ws.Cells.Item(3,15).Value = 0.01900486009268527
$ws.Cells.Item(3,16).Value = 0.01900486009268527
$ws.Cells.Item(3,17).Value = 15.26279925244467
$ws.Cells.Item(3,18).Value = 137.365193272002
$ws.Cells.Item(3,19).Value = 0.00373129353890596
$ws.Cells.Item(3,20).Value = 0.00373129353890596
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Cd9"
$ws.Cells.Item(4,3).Value = "L1cam"
$ws.Cells.Item(4,4).Value = "M2"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 40.314886
$ws.Cells.Item(4,8).Value = 120.944658
$ws.Cells.Item(4,9).Value = 0.1963336494301312
$ws.Cells.Item(4,10).Value = 0.1963336494301312
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 4.123197999999999
$ws.Cells.Item(4,14).Value = 12.369594
$ws.Cells.Item(4,15).Value = 0.2069808238940481
$ws.Cells.Item(4,16).Value = 0.2069808238940482
$ws.Cells.Item(4,17).Value = 166.226257325428
$ws.Cells.Item(4,18).Value = 1496.036315928852
$ws.Cells.Item(4,19).Value = 0.04063730051717376
$ws.Cells.Item(4,20).Value = 0.04063730051717377
$ws.Cells.Item(5,1).Value = "ECs"
$ws.Cells.Item(5,2).Value = "Cd9"
$ws.Cells.Item(5,3).Value = "L1cam"
$ws.Cells.Item(5,4).Value = "sCs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 40.314886
$ws.Cells.Item(5,8).Value = 120.944658
$ws.Cells.Item(5,9).Value = 0.1963336494301312
$ws.Cells.Item(5,10).Value = 0.1963336494301312
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 3.531764333333333
$ws.Cells.Item(5,14).Value = 10.595293
$ws.Cells.Item(5,15).Value = 0.1772913868101768
$ws.Cells.Item(5,16).Value = 0.1772913868101768
$ws.Cells.Item(5,17).Value = 142.3826764771993
$ws.Cells.Item(5,18).Value = 1281.444088294794
$ws.Cells.Item(5,19).Value = 0.03480826498497102
$ws.Cells.Item(5,20).Value = 0.03480826498497102
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Cd9"
$ws.Cells.Item(6,3).Value = "L1cam"
$ws.Cells.Item(6,4).Value = "ECs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 69.65329366666667
$ws.Cells.Item(6,8).Value = 208.959881
$ws.Cells.Item(6,9).Value = 0.3392118072814421
$ws.Cells.Item(6,10).Value = 0.3392118072814421
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 11.88712433333333
$ws.Cells.Item(6,14).Value = 35.661373
$ws.Cells.Item(6,15).Value = 0.5967229292030898
$ws.Cells.Item(6,16).Value = 0.5967229292030898
$ws.Cells.Item(6,17).Value = 827.9773620418458
$ws.Cells.Item(6,18).Value = 7451.796258376612
$ws.Cells.Item(6,19).Value = 0.2024154632612561
$ws.Cells.Item(6,20).Value = 0.2024154632612561
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Cd9"
$ws.Cells.Item(7,3).Value = "L1cam"
$ws.Cells.Item(7,4).Value = "FAPs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 69.65329366666667
$ws.Cells.Item(7,8).Value = 208.959881
$ws.Cells.Item(7,9).Value = 0.3392118072814421
$ws.Cells.Item(7,10).Value = 0.3392118072814421
$ws.Cells.Item(7,11).Value = 2
$ws.Cells.Item(7,12).Value = 0.6666666666666666
$ws.Cells.Item(7,13).Value = 0.3785896666666667
$ws.Cells.Item(7,14).Value = 1.135769
$ws.Cells.Item(7,15).Value = 0.01900486009268527
$ws.Cells.Item(7,16).Value = 0.01900486009268527
$ws.Cells.Item(7,17).Value = 26.37001723149878
$ws.Cells.Item(7,18).Value = 237.330155083489
$ws.Cells.Item(7,19).Value = 0.006446672939170727
$ws.Cells.Item(7,20).Value = 0.006446672939170727
$ws.Cells.Item(8,1).Value = "FAPs"
$ws.Cells.Item(8,2).Value = "Cd9"
$ws.Cells.Item(8,3).Value = "L1cam"
$ws.Cells.Item(8,4).Value = "M2"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 69.65329366666667
$ws.Cells.Item(8,8).Value = 208.959881
$ws.Cells.Item(8,9).Value = 0.3392118072814421
$ws.Cells.Item(8,10).Value = 0.3392118072814421
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 4.123197999999999
$ws.Cells.Item(8,14).Value = 12.369594
$ws.Cells.Item(8,15).Value = 0.2069808238940481
$ws.Cells.Item(8,16).Value = 0.2069808238940482
$ws.Cells.Item(8,17).Value = 287.1943211398126
$ws.Cells.Item(8,18).Value = 2584.748890258314
$ws.Cells.Item(8,19).Value = 0.07021033934570196
$ws.Cells.Item(8,20).Value = 0.07021033934570196
$ws.Cells.Item(9,1).Value = "FAPs"
$ws.Cells.Item(9,2).Value = "Cd9"
$ws.Cells.Item(9,3).Value = "L1cam"
$ws.Cells.Item(9,4).Value = "sCs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 69.65329366666667
$ws.Cells.Item(9,8).Value = 208.959881
$ws.Cells.Item(9,9).Value = 0.3392118072814421
$ws.Cells.Item(9,10).Value = 0.3392118072814421
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 3.531764333333333
$ws.Cells.Item(9,14).Value = 10.595293
$ws.Cells.Item(9,15).Value = 0.1772913868101768
$ws.Cells.Item(9,16).Value = 0.1772913868101768
$ws.Cells.Item(9,17).Value = 245.9990182711259
$ws.Cells.Item(9,18).Value = 2213.991164440133
$ws.Cells.Item(9,19).Value = 0.06013933173531328
$ws.Cells.Item(9,20).Value = 0.06013933173531328
$ws.Cells.Item(10,1).Value = "M2"
$ws.Cells.Item(10,2).Value = "Cd9"
$ws.Cells.Item(10,3).Value = "L1cam"
$ws.Cells.Item(10,4).Value = "ECs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 69.40355333333332
$ws.Cells.Item(10,8).Value = 208.21066
$ws.Cells.Item(10,9).Value = 0.3379955709003387
$ws.Cells.Item(10,10).Value = 0.3379955709003388
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 11.88712433333333
$ws.Cells.Item(10,14).Value = 35.661373
$ws.Cells.Item(10,15).Value = 0.5967229292030898
$ws.Cells.Item(10,16).Value = 0.5967229292030898
$ws.Cells.Item(10,17).Value = 825.0086676484642
$ws.Cells.Item(10,18).Value = 7425.078008836178
$ws.Cells.Item(10,19).Value = 0.2016897071253207
$ws.Cells.Item(10,20).Value = 0.2016897071253208
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Cd9"
$ws.Cells.Item(11,3).Value = "L1cam"
$ws.Cells.Item(11,4).Value = "FAPs"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 69.40355333333332
$ws.Cells.Item(11,8).Value = 208.21066
$ws.Cells.Item(11,9).Value = 0.3379955709003387
$ws.Cells.Item(11,10).Value = 0.3379955709003388
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.3785896666666667
$ws.Cells.Item(11,14).Value = 1.135769
$ws.Cells.Item(11,15).Value = 0.01900486009268527
$ws.Cells.Item(11,16).Value = 0.01900486009268527
$ws.Cells.Item(11,17).Value = 26.27546812194888
$ws.Cells.Item(11,18).Value = 236.47921309754
$ws.Cells.Item(11,19).Value = 0.006423558536908224
$ws.Cells.Item(11,20).Value = 0.006423558536908225
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Cd9"
$ws.Cells.Item(12,3).Value = "L1cam"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 69.40355333333332
$ws.Cells.Item(12,8).Value = 208.21066
$ws.Cells.Item(12,9).Value = 0.3379955709003387
$ws.Cells.Item(12,10).Value = 0.3379955709003388
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 4.123197999999999
$ws.Cells.Item(12,14).Value = 12.369594
$ws.Cells.Item(12,15).Value = 0.2069808238940481
$ws.Cells.Item(12,16).Value = 0.2069808238940482
$ws.Cells.Item(12,17).Value = 286.1645922968933
$ws.Cells.Item(12,18).Value = 2575.481330672039
$ws.Cells.Item(12,19).Value = 0.06995860173749127
$ws.Cells.Item(12,20).Value = 0.0699586017374913
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Cd9"
$ws.Cells.Item(13,3).Value = "L1cam"
$ws.Cells.Item(13,4).Value = "sCs"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 69.40355333333332
$ws.Cells.Item(13,8).Value = 208.21066
$ws.Cells.Item(13,9).Value = 0.3379955709003387
$ws.Cells.Item(13,10).Value = 0.3379955709003388
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 3.531764333333333
$ws.Cells.Item(13,14).Value = 10.595293
$ws.Cells.Item(13,15).Value = 0.1772913868101768
$ws.Cells.Item(13,16).Value = 0.1772913868101768
$ws.Cells.Item(13,17).Value = 245.1169942692644
$ws.Cells.Item(13,18).Value = 2206.052948423379
$ws.Cells.Item(13,19).Value = 0.05992370350061847
$ws.Cells.Item(13,20).Value = 0.05992370350061849
$ws.Cells.Item(14,1).Value = "sCs"
$ws.Cells.Item(14,2).Value = "Cd9"
$ws.Cells.Item(14,3).Value = "L1cam"
$ws.Cells.Item(14,4).Value = "ECs"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 25.96691433333334
$ws.Cells.Item(14,8).Value = 77.90074300000001
$ws.Cells.Item(14,9).Value = 0.126458972388088
$ws.Cells.Item(14,10).Value = 0.126458972388088
$ws.Cells.Item(14,11).Value = 3
$ws.Cells.Item(14,12).Value = 1
$ws.Cells.Item(14,13).Value = 11.88712433333333
$ws.Cells.Item(14,14).Value = 35.661373
$ws.Cells.Item(14,15).Value = 0.5967229292030898
$ws.Cells.Item(14,16).Value = 0.5967229292030898
$ws.Cells.Item(14,17).Value = 308.6719392333488
$ws.Cells.Item(14,18).Value = 2778.047453100139
$ws.Cells.Item(14,19).Value = 0.07546096842743251
$ws.Cells.Item(14,20).Value = 0.07546096842743251
$ws.Cells.Item(15,1).Value = "sCs"
$ws.Cells.Item(15,2).Value = "Cd9"
$ws.Cells.Item(15,3).Value = "L1cam"
$ws.Cells.Item(15,4).Value = "FAPs"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 25.96691433333334
$ws.Cells.Item(15,8).Value = 77.90074300000001
$ws.Cells.Item(15,9).Value = 0.126458972388088
$ws.Cells.Item(15,10).Value = 0.126458972388088
$ws.Cells.Item(15,11).Value = 2
$ws.Cells.Item(15,12).Value = 0.6666666666666666
$ws.Cells.Item(15,13).Value = 0.3785896666666667
$ws.Cells.Item(15,14).Value = 1.135769
$ws.Cells.Item(15,15).Value = 0.01900486009268527
$ws.Cells.Item(15,16).Value = 0.01900486009268527
$ws.Cells.Item(15,17).Value = 9.830805441818557
$ws.Cells.Item(15,18).Value = 88.47724897636701
$ws.Cells.Item(15,19).Value = 0.002403335077700362
$ws.Cells.Item(15,20).Value = 0.002403335077700362
$ws.Cells.Item(16,1).Value = "sCs"
$ws.Cells.Item(16,2).Value = "Cd9"
$ws.Cells.Item(16,3).Value = "L1cam"
$ws.Cells.Item(16,4).Value = "M2"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 25.96691433333334
$ws.Cells.Item(16,8).Value = 77.90074300000001
$ws.Cells.Item(16,9).Value = 0.126458972388088
$ws.Cells.Item(16,10).Value = 0.126458972388088
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 4.123197999999999
$ws.Cells.Item(16,14).Value = 12.369594
$ws.Cells.Item(16,15).Value = 0.2069808238940481
$ws.Cells.Item(16,16).Value = 0.2069808238940482
$ws.Cells.Item(16,17).Value = 107.0667292453713
$ws.Cells.Item(16,18).Value = 963.600563208342
$ws.Cells.Item(16,19).Value = 0.02617458229368114
$ws.Cells.Item(16,20).Value = 0.02617458229368114
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Cd9"
$ws.Cells.Item(17,3).Value = "L1cam"
$ws.Cells.Item(17,4).Value = "sCs"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 25.96691433333334
$ws.Cells.Item(17,8).Value = 77.90074300000001
$ws.Cells.Item(17,9).Value = 0.126458972388088
$ws.Cells.Item(17,10).Value = 0.126458972388088
$ws.Cells.Item(17,11).Value = 3
$ws.Cells.Item(17,12).Value = 1
$ws.Cells.Item(17,13).Value = 3.531764333333333
$ws.Cells.Item(17,14).Value = 10.595293
$ws.Cells.Item(17,15).Value = 0.1772913868101768
$ws.Cells.Item(17,16).Value = 0.1772913868101768
$ws.Cells.Item(17,17).Value = 91.70902188918879
$ws.Cells.Item(17,18).Value = 825.3811970026991
$ws.Cells.Item(17,19).Value = 0.02242008658927397
$ws.Cells.Item(17,20).Value = 0.02242008658927397